$d = $word.ActiveDocument

$d.Content.Find.Execute("66÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "28÷2=", 2) | Out-Null
$d.Content.Find.Execute("10÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "12÷8=", 2) | Out-Null
$d.Content.Find.Execute("48÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "78÷2=", 2) | Out-Null
$d.Content.Find.Execute("28÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "47÷7=", 2) | Out-Null
$d.Content.Find.Execute("82÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "65÷5=", 2) | Out-Null
$d.Content.Find.Execute("84÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "21÷2=", 2) | Out-Null
$d.Content.Find.Execute("85÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "21÷7=", 2) | Out-Null
$d.Content.Find.Execute("70÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "36÷7=", 2) | Out-Null
$d.Content.Find.Execute("94÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "80÷3=", 2) | Out-Null
$d.Content.Find.Execute("94÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "92÷2=", 2) | Out-Null
$d.Content.Find.Execute("83÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "45÷6=", 2) | Out-Null
$d.Content.Find.Execute("78÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "85÷2=", 2) | Out-Null
$d.Content.Find.Execute("40÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "35÷4=", 2) | Out-Null
$d.Content.Find.Execute("48÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "61÷8=", 2) | Out-Null
$d.Content.Find.Execute("14÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "47÷3=", 2) | Out-Null
$d.Content.Find.Execute("30÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "71÷7=", 2) | Out-Null
$d.Content.Find.Execute("96÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "88÷3=", 2) | Out-Null
$d.Content.Find.Execute("31÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "70÷6=", 2) | Out-Null
$d.Content.Find.Execute("76÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "18÷2=", 2) | Out-Null
$d.Content.Find.Execute("11÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "91÷8=", 2) | Out-Null
$d.Content.Find.Execute("40÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "68÷5=", 2) | Out-Null
$d.Content.Find.Execute("76÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "71÷7=", 2) | Out-Null
$d.Content.Find.Execute("81÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "37÷2=", 2) | Out-Null
$d.Content.Find.Execute("24÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "43÷2=", 2) | Out-Null
